$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing description for Atolla wyvillei (row 4) ---
$ws.Range("E4").Value = "Questa specie di medusa presenta una campana rossa e un lungo tentacolo bianco che parte dal centro della campana e si allunga per una decina di centimetri"

# --- Fill in new organism rows (6-10) ---

# Row 6: Marrus orthocana / Medusa sifonofora cintura di fuoco
$ws.Range("B6").Value = "Marrus orthocana"
$ws.Range("C6").Value = "null"
$ws.Range("D6").Value = "Medusa sifonofora cintura di fuoco"
$ws.Range("E6").Value = "Questa specie di sifonoforo presenta un colore rosso/arancio intenso, essendo i sifonofori dei raggruppamenti di organismi, può assumere diverse forme, caratterizzate solitamente da filamenti intrecciati."

# Row 7: Solmissus / Medusa piatto da cena
$ws.Range("B7").Value = "Solmissus"
$ws.Range("C7").Value = "null"
$ws.Range("D7").Value = "Medusa piatto da cena"
$ws.Range("E7").Value = "Questa specie di medusa presenta numerosi e sottili tentacoli. La campana è quasi completamente trasparente. Sia i tentacoli che la campana sono fluorescenti e di colore bianco/azzurro."

# Row 8: Tiburonia granrojo / Grande medusa rossa
$ws.Range("B8").Value = "Tiburonia granrojo"
$ws.Range("C8").Value = "null"
$ws.Range("D8").Value = "Grande medusa rossa"
$ws.Range("E8").Value = "Il nome di questa medusa deriva dai suoi colori rossi della campana (granrojo = grande rosso in spagnolo). Presenta una campana di grandi dimensioni, mentre i tentacoli sono più corti e tozzi."

# Row 9: Medusa quadri-tentacolo
$ws.Range("B9").Value = "null"
$ws.Range("C9").Value = "Medusa quadri-tentacolo"
$ws.Range("D9").Value = "null"
$ws.Range("E9").Value = "Questa medusa di piccole dimensioni presenta 4 sottili tentacoli e una campana trasparente con un anello fluorescente. Utilizza la propulsione causata dal movimento della campana per muoversi."

# Row 10: Larvacean
$ws.Range("B10").Value = "Larvacean"
$ws.Range("C10").Value = "null"
$ws.Range("D10").Value = "null"
$ws.Range("E10").Value = "I larvacei sono tunicati che abitano diverse zone degli oceani di tutto il mondo, per nutrirsi utilizzano una “barriera” di muco che filtra le sostanze nutritive."

# --- Adjust row heights to fit the new wrapped text ---
$ws.Rows.Item(4).RowHeight = 27.7
$ws.Rows.Item(6).RowHeight = 40.95
$ws.Rows.Item(7).RowHeight = 27.7
$ws.Rows.Item(8).RowHeight = 40.95
$ws.Rows.Item(9).RowHeight = 27.7
$ws.Rows.Item(10).RowHeight = 27.7

# --- Update the active selection to reflect where editing ended ---
$null = $ws.Range("E12").Select()
